$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 15.05048780176053
$ws.Range("C2").Value = 11.50101368500164
$ws.Range("D2").Value = 3.627452838550872
$ws.Range("E2").Value = 16.59448567039335
$ws.Range("F2").Value = 19.02170427638955
$ws.Range("I2").Value = 15.77803083377724
$ws.Range("O2").Value = 16.8512498480365
$ws.Range("B3").Value = 14.20706712963875
$ws.Range("C3").Value = 10.82536328834251
$ws.Range("D3").Value = 3.574240261778049
$ws.Range("E3").Value = 15.64601514548837
$ws.Range("F3").Value = 19.02300677213631
$ws.Range("I3").Value = 15.94719198431599
$ws.Range("O3").Value = 16.93076996332746
$ws.Range("B4").Value = 13.66225859234334
$ws.Range("C4").Value = 10.38656213626149
$ws.Range("D4").Value = 3.541051168541962
$ws.Range("E4").Value = 15.038347417373
$ws.Range("F4").Value = 19.03334216288247
$ws.Range("I4").Value = 16.05675514797745
$ws.Range("O4").Value = 16.98693378145153
$ws.Range("B5").Value = 13.43362357182914
$ws.Range("C5").Value = 10.20178148734509
$ws.Range("D5").Value = 3.527407687069452
$ws.Range("E5").Value = 14.78462218534459
$ws.Range("F5").Value = 19.03993903451894
$ws.Range("I5").Value = 16.10283464319667
$ws.Range("O5").Value = 17.011649856549
$ws.Range("B6").Value = 13.39526417595278
$ws.Range("C6").Value = 10.17074045402317
$ws.Range("D6").Value = 3.525135363012186
$ws.Range("E6").Value = 14.7421319904535
$ws.Range("F6").Value = 19.04117806589743
$ws.Range("I6").Value = 16.11057254693057
$ws.Range("O6").Value = 17.01586391232873
$ws.Range("B7").Value = 13.65920172075717
$ws.Range("C7").Value = 10.38409418832189
$ws.Range("D7").Value = 3.540867633260101
$ws.Range("E7").Value = 15.03494987954568
$ws.Range("F7").Value = 19.03342149243547
$ws.Range("I7").Value = 16.05737079695219
$ws.Range("O7").Value = 16.98725972819071
$ws.Range("B8").Value = 14.76535150714479
$ws.Range("C8").Value = 11.27304942293409
$ws.Range("D8").Value = 3.609218932299791
$ws.Range("E8").Value = 16.27284911097061
$ws.Range("F8").Value = 19.02016752969193
$ws.Range("I8").Value = 15.83517438231825
$ws.Range("O8").Value = 16.87713566352617
$ws.Range("B9").Value = 16.71544598605509
$ws.Range("C9").Value = 12.82448042983457
$ws.Range("D9").Value = 3.738683220059439
$ws.Range("E9").Value = 18.60857781505257
$ws.Range("F9").Value = 19.07025629307025
$ws.Range("I9").Value = 15.44470843078692
$ws.Range("O9").Value = 16.72011323027554
$ws.Range("B10").Value = 18.00929717279027
$ws.Range("C10").Value = 13.84578044727016
$ws.Range("D10").Value = 3.830401932369729
$ws.Range("E10").Value = 20.26597166812811
$ws.Range("F10").Value = 19.15378772040472
$ws.Range("I10").Value = 15.18550722359725
$ws.Range("O10").Value = 16.64159668072431
$ws.Range("B11").Value = 18.56700732572864
$ws.Range("C11").Value = 14.28448812996355
$ws.Range("D11").Value = 3.871267465152992
$ws.Range("E11").Value = 20.97755886320545
$ws.Range("F11").Value = 19.20194077485662
$ws.Range("I11").Value = 15.07362239834789
$ws.Range("O11").Value = 16.61406113290711
$ws.Range("B12").Value = 18.77371547950815
$ws.Range("C12").Value = 14.44688581814688
$ws.Range("D12").Value = 3.886609214091045
$ws.Range("E12").Value = 21.24096322204026
$ws.Range("F12").Value = 19.22163034460086
$ws.Range("I12").Value = 15.03212393127967
$ws.Range("O12").Value = 16.60482444454607
$ws.Range("B13").Value = 18.72939723691359
$ws.Range("C13").Value = 14.41207661595647
$ws.Range("D13").Value = 3.883311166422821
$ws.Range("E13").Value = 21.18450327772759
$ws.Range("F13").Value = 19.21732524687685
$ws.Range("I13").Value = 15.04102262547547
$ws.Range("O13").Value = 16.6067605752739
$ws.Range("B14").Value = 18.5841034277223
$ws.Range("C14").Value = 14.2979235020704
$ws.Range("D14").Value = 3.872532362862923
$ws.Range("E14").Value = 20.99935051864087
$ws.Range("F14").Value = 19.20353152746721
$ws.Range("I14").Value = 15.0701908410748
$ws.Range("O14").Value = 16.61327729098467
$ws.Range("B15").Value = 18.49452159502535
$ws.Range("C15").Value = 14.22751530486794
$ws.Range("D15").Value = 3.865912412793098
$ws.Range("E15").Value = 20.88515116012722
$ws.Range("F15").Value = 19.19527175751346
$ws.Range("I15").Value = 15.08817060131386
$ws.Range("O15").Value = 16.61742439175776
$ws.Range("B16").Value = 17.97222416279133
$ws.Range("C16").Value = 13.81658808603135
$ws.Range("D16").Value = 3.82771316521258
$ws.Range("E16").Value = 20.21861813496901
$ws.Range("F16").Value = 19.15084477423367
$ws.Range("I16").Value = 15.19294070623337
$ws.Range("O16").Value = 16.64356204246265
$ws.Range("B17").Value = 17.64386893414432
$ws.Range("C17").Value = 13.55785937866145
$ws.Range("D17").Value = 3.804052352086922
$ws.Range("E17").Value = 19.79889441212995
$ws.Range("F17").Value = 19.12618843560828
$ws.Range("I17").Value = 15.25875962396454
$ws.Range("O17").Value = 16.66170240610372
$ws.Range("B18").Value = 17.45210401644001
$ws.Range("C18").Value = 13.40660905857839
$ws.Range("D18").Value = 3.790363097650111
$ws.Range("E18").Value = 19.55349035728542
$ws.Range("F18").Value = 19.11296304517519
$ws.Range("I18").Value = 15.29718401207258
$ws.Range("O18").Value = 16.67290562093229
$ws.Range("B19").Value = 17.3866785885234
$ws.Range("C19").Value = 13.35498019618174
$ws.Range("D19").Value = 3.785714686303463
$ws.Range("E19").Value = 19.4697144103139
$ws.Range("F19").Value = 19.10864949766634
$ws.Range("I19").Value = 15.31029115561299
$ws.Range("O19").Value = 16.67683053756916
$ws.Range("B20").Value = 17.67912379592172
$ws.Range("C20").Value = 13.58565373029128
$ws.Range("D20").Value = 3.806579456878181
$ws.Range("E20").Value = 19.84398742526866
$ws.Range("F20").Value = 19.12871419968848
$ws.Range("I20").Value = 15.25169437456008
$ws.Range("O20").Value = 16.65969160476427
$ws.Range("B21").Value = 18.62690175775741
$ws.Range("C21").Value = 14.33155435474239
$ws.Range("D21").Value = 3.875702047906167
$ws.Range("E21").Value = 21.05389851854794
$ws.Range("F21").Value = 19.20754364497656
$ws.Range("I21").Value = 15.06159979336352
$ws.Range("O21").Value = 16.61133075990861
$ws.Range("B22").Value = 19.22017648461896
$ws.Range("C22").Value = 14.79729320260919
$ws.Range("D22").Value = 3.920097397021088
$ws.Range("E22").Value = 21.80934895646643
$ws.Range("F22").Value = 19.26753966070871
$ws.Range("I22").Value = 14.94243410607024
$ws.Range("O22").Value = 16.58666862874338
$ws.Range("B23").Value = 18.90593976499226
$ws.Range("C23").Value = 14.55071151239328
$ws.Range("D23").Value = 3.896477271721762
$ws.Range("E23").Value = 21.40936914539612
$ws.Range("F23").Value = 19.23474562453683
$ws.Range("I23").Value = 15.00556982922822
$ws.Range("O23").Value = 16.59919151526063
$ws.Range("B24").Value = 17.66319437979707
$ws.Range("C24").Value = 13.57309570185001
$ws.Range("D24").Value = 3.805437220963207
$ws.Range("E24").Value = 19.82361364538989
$ws.Range("F24").Value = 19.12756934266482
$ws.Range("I24").Value = 15.25488675424893
$ws.Range("O24").Value = 16.66059827792072
$ws.Range("B25").Value = 16.21197793584842
$ws.Range("C25").Value = 12.42551268275503
$ws.Range("D25").Value = 3.704213794196838
$ws.Range("E25").Value = 17.96042625348758
$ws.Range("F25").Value = 19.048511898994
$ws.Range("I25").Value = 15.54548694947148
$ws.Range("O25").Value = 16.75618514696082
